$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Characters(21, 2).Text = "46"
# Replace the rightmost (2nd) date first so the leftmost date offset
# is unaffected by the length change (9 chars -> 11 chars).
$ws.Range("C9").Characters(47, 10).Text = "11/17/2024"
$ws.Range("C9").Characters(27, 9).Text = "11/11/2024"

# --- Crime statistics table updates (rows 14-33) ---
# Cells whose type/style changes: copy a donor cell (matching target
# style+type) first, then (for numeric targets) overwrite with the value.
$ws.Range("C14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("I14").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 3
$ws.Range("K14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = 33.333333333333
$ws.Range("I14").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("K14").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("I14").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("K14").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100
$ws.Range("I14").Copy($ws.Range("C33"))
$ws.Range("C33").Value = 2

# Cells with same style, value-only updates
$ws.Range("M14").Value = -15.384615384615
$ws.Range("F15").Value = 6
$ws.Range("H15").Value = 100
$ws.Range("L15").Value = -33.333333333333
$ws.Range("N15").Value = -54.545454545454
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 14
$ws.Range("E16").Value = -78.571428571428
$ws.Range("F16").Value = 39
$ws.Range("G16").Value = 43
$ws.Range("H16").Value = -9.302325581395
$ws.Range("I16").Value = 438
$ws.Range("J16").Value = 476
$ws.Range("K16").Value = -7.98319327731
$ws.Range("L16").Value = -33.93665158371
$ws.Range("M16").Value = -4.366812227074
$ws.Range("N16").Value = -75.973669775096
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = 14.285714285714
$ws.Range("F17").Value = 61
$ws.Range("H17").Value = 24.489795918367
$ws.Range("I17").Value = 688
$ws.Range("J17").Value = 748
$ws.Range("K17").Value = -8.021390374331
$ws.Range("L17").Value = -2.549575070821
$ws.Range("M17").Value = 52.212389380531
$ws.Range("N17").Value = -18.483412322274
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = -30
$ws.Range("I18").Value = 273
$ws.Range("J18").Value = 285
$ws.Range("K18").Value = -4.210526315789
$ws.Range("L18").Value = -1.444043321299
$ws.Range("M18").Value = -20.408163265306
$ws.Range("N18").Value = -83.424408014572
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -5.263157894736
$ws.Range("F19").Value = 86
$ws.Range("G19").Value = 72
$ws.Range("H19").Value = 19.444444444444
$ws.Range("I19").Value = 929
$ws.Range("J19").Value = 818
$ws.Range("K19").Value = 13.569682151589
$ws.Range("L19").Value = 1.19825708061
$ws.Range("M19").Value = 83.596837944664
$ws.Range("N19").Value = 33.477011494252
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 14
$ws.Range("E20").Value = -42.857142857142
$ws.Range("F20").Value = 50
$ws.Range("G20").Value = 45
$ws.Range("H20").Value = 11.111111111111
$ws.Range("I20").Value = 479
$ws.Range("J20").Value = 576
$ws.Range("K20").Value = -16.840277777777
$ws.Range("L20").Value = 1.914893617021
$ws.Range("M20").Value = 120.73732718894
$ws.Range("N20").Value = -72.998872604284
$ws.Range("C21").Value = 49
$ws.Range("D21").Value = 67
$ws.Range("E21").Value = -26.865671641791
$ws.Range("F21").Value = 263
$ws.Range("G21").Value = 242
$ws.Range("H21").Value = 8.677685950413
$ws.Range("I21").Value = 2848
$ws.Range("J21").Value = 2953
$ws.Range("K21").Value = -3.555706061632
$ws.Range("L21").Value = -7.682333873581
$ws.Range("M21").Value = 40.850642927794
$ws.Range("N21").Value = -58.790334249746
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 300
$ws.Range("L22").Value = 30
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 400
$ws.Range("F23").Value = 31
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = 82.35294117647
$ws.Range("I23").Value = 261
$ws.Range("J23").Value = 275
$ws.Range("K23").Value = -5.090909090909
$ws.Range("L23").Value = -11.824324324324
$ws.Range("M23").Value = 31.155778894472
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -8.108108108108
$ws.Range("F24").Value = 138
$ws.Range("G24").Value = 142
$ws.Range("H24").Value = -2.81690140845
$ws.Range("I24").Value = 1660
$ws.Range("J24").Value = 1721
$ws.Range("K24").Value = -3.544450900639
$ws.Range("L24").Value = -10.124526258798
$ws.Range("M24").Value = 25.472411186696
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -9.090909090909
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = -11.111111111111
$ws.Range("I25").Value = 645
$ws.Range("J25").Value = 703
$ws.Range("K25").Value = -8.250355618776
$ws.Range("L25").Value = -36.390532544378
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 22
$ws.Range("E26").Value = -36.363636363636
$ws.Range("G26").Value = 71
$ws.Range("H26").Value = 8.450704225352
$ws.Range("I26").Value = 1042
$ws.Range("J26").Value = 981
$ws.Range("K26").Value = 6.218144750254
$ws.Range("L26").Value = 3.992015968063
$ws.Range("M26").Value = -24.437998549673
$ws.Range("F27").Value = 8
$ws.Range("H27").Value = 60
$ws.Range("L27").Value = -28.787878787878
$ws.Range("C28").Value = 4
$ws.Range("F28").Value = 10
$ws.Range("H28").Value = 66.666666666666
$ws.Range("I28").Value = 101
$ws.Range("J28").Value = 98
$ws.Range("K28").Value = 3.061224489795
$ws.Range("L28").Value = 60.31746031746
$ws.Range("J29").Value = 39
$ws.Range("K29").Value = -33.333333333333
$ws.Range("L29").Value = -36.585365853658
$ws.Range("M29").Value = -43.478260869565
$ws.Range("N29").Value = -81.294964028777
$ws.Range("J30").Value = 33
$ws.Range("K30").Value = -30.30303030303
$ws.Range("L30").Value = -32.35294117647
$ws.Range("M30").Value = -41.025641025641
$ws.Range("N30").Value = -81.6
$ws.Range("F33").Value = 3
$ws.Range("I33").Value = 9
$ws.Range("K33").Value = 50
$ws.Range("L33").Value = 0
